$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.305.65'
$ws.Range('E2').Value = '  +1.30%  '
$ws.Range('D3').Value = '1.624.43'
$ws.Range('E3').Value = '  +1.63%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '''212.71'
$ws.Range('E5').Value = '  +0.75%  '
$ws.Range('E6').Value = '  -0.03%  '
$ws.Range('E7').Value = '  +0.48%  '
$ws.Range('D8').Value = '''0.250'
$ws.Range('E8').Value = '  +1.81%  '
$ws.Range('D9').Value = '''0.0616'
$ws.Range('E9').Value = '  +0.75%  '
$ws.Range('D10').Value = '''19.00'
$ws.Range('E10').Value = '  +5.15%  '
$ws.Range('D11').Value = '''0.0815'
$ws.Range('E11').Value = '  +0.40%  '
$ws.Range('D12').Value = '1.849.33'
$ws.Range('E12').Value = '  +1.49%  '
$ws.Range('D13').Value = '1.648.67'
$ws.Range('E13').Value = '  +3.07%  '
$ws.Range('E14').Value = '  +0.79%  '
$ws.Range('D15').Value = '''0.520'
$ws.Range('E15').Value = '  +1.33%  '
$ws.Range('D16').Value = '26.321.14'
$ws.Range('E16').Value = '  +1.27%  '
$ws.Range('D17').Value = '''62.49'
$ws.Range('E17').Value = '  +3.90%  '
$ws.Range('D18').Value = '0.0₃0730'
$ws.Range('E18').Value = '  +1.11%  '
$ws.Range('E19').Value = '  -0.02%  '
$ws.Range('D20').Value = '''202.99'
$ws.Range('E20').Value = '  +0.83%  '
$ws.Range('E21').Value = '  +1.75%  '
$ws.Range('E22').Value = '  +1.39%  '
$ws.Range('D23').Value = '''6.06'
$ws.Range('E23').Value = '  +0.95%  '
$ws.Range('D24').Value = '''1.92'
$ws.Range('E24').Value = '  +7.17%  '
$ws.Range('D25').Value = '''142.78'
$ws.Range('E25').Value = '  +0.89%  '
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('E27').Value = '  -0.21%  '
$ws.Range('E28').Value = '  +0.87%  '
$ws.Range('E29').Value = '  +1.82%  '
$ws.Range('D30').Value = '''0.0526'
$ws.Range('E30').Value = '  +10.60%  '
$ws.Range('E31').Value = '  +0.59%  '
$ws.Range('D32').Value = '''3.19'
$ws.Range('E32').Value = '  +2.88%  '
$ws.Range('D33').Value = '''2.96'
$ws.Range('E33').Value = '  +0.15%  '
$ws.Range('E34').Value = '  +2.36%  '
$ws.Range('E35').Value = '  +3.06%  '
$ws.Range('D36').Value = '1.177.34'
$ws.Range('E36').Value = '  +4.86%  '
$ws.Range('E37').Value = '  +1.19%  '
$ws.Range('E38').Value = '  +3.34%  '
$ws.Range('E39').Value = '  -0.01%  '
$ws.Range('E40').Value = '  -0.15%  '
$ws.Range('E41').Value = '  +1.73%  '
$ws.Range('D42').Value = '''0.794'
$ws.Range('E42').Value = '  +1.06%  '
$ws.Range('D43').Value = '''5.33'
$ws.Range('E43').Value = '  +3.82%  '
$ws.Range('D44').Value = '1.760.48'
$ws.Range('E44').Value = '  +1.53%  '
$ws.Range('D45').Value = '''93.49'
$ws.Range('E45').Value = '  +0.79%  '
$ws.Range('E46').Value = '  +14.90%  '
$ws.Range('E47').Value = '  +1.38%  '
$ws.Range('D48').Value = '''54.18'
$ws.Range('E48').Value = '  +1.56%  '
$ws.Range('E49').Value = '  +1.03%  '
$ws.Range('D50').Value = '''0.408'
$ws.Range('E50').Value = '  +0.08%  '
$ws.Range('E51').Value = '  -0.27%  '
